# Remove stray trailing spaces left on several lyric lines (trailing
# whitespace before the line break) across slides 1, 2, 3 and 4.
#
# Each affected shape is an auto-fit text box ("TextBox 1"); nudging its
# TextRange recalculates the box's rendered height, so we snapshot and
# restore Height on every shape we touch to avoid an unrelated side effect.

$p = $ppt.ActivePresentation

function Remove-TrailingSpace($shape, $pos) {
    # Deletes the single character at 1-based position $pos (expected to
    # be the offending trailing space) by collapsing it to an empty run.
    $shape.TextFrame.TextRange.Characters($pos, 1).Text = ""
}

# --- Slide 1: "Great is Your faithfulness Oh God ..." ---
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(1)
$h1 = $shp1.Height
Remove-TrailingSpace $shp1 142   # "And nothing can keep us apart "
Remove-TrailingSpace $shp1 72    # "You wrestle with the sinner's heart "
Remove-TrailingSpace $shp1 35    # "Great is Your faithfulness Oh God "
$shp1.Height = $h1

# --- Slide 2: "Remember Your children / Remember Your promise Oh God" ---
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(1)
$h2 = $shp2.Height
Remove-TrailingSpace $shp2 78    # "Remember Your promise Oh God "
Remove-TrailingSpace $shp2 48    # "Remember Your children "
$shp2.Height = $h2

# --- Slide 3: "Your grace is enough" (first line only) ---
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(1)
$h3 = $shp3.Height
Remove-TrailingSpace $shp3 22    # "Your grace is enough "
$shp3.Height = $h3

# --- Slide 4: "Great is Your love and justice God ..." ---
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(1)
$h4 = $shp4.Height
Remove-TrailingSpace $shp4 147   # "And all Your people sing along "
Remove-TrailingSpace $shp4 73    # "You use the weak to lead the strong "
Remove-TrailingSpace $shp4 36    # "Great is Your love and justice God "
$shp4.Height = $h4
